$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged, values updated)
$ws.Range("B3").Value = 0.9988667617076624
$ws.Range("C3").Value = 0.9988399599154901
$ws.Range("D3").Value = 0.9904959739028952

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9987617446009914
$ws.Range("C4").Value = 0.9985969538182865
$ws.Range("D4").Value = 0.9796081043219118

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9992949010265352
$ws.Range("C5").Value = 0.9990730538946045
$ws.Range("D5").Value = 0.998253275466078
